{"js": "// Replace the three-digit x one-digit multiplication problems/answers\n// in the table cells with the new values from the commit diff.\n// Each old value is unique within the document, so a plain text search\n// (whole document body) followed by a full-text replace of the matched\n// range is sufficient and keeps the original run formatting intact.\n\nconst replacements = [\n  [\"276\u00d77=1932\", \"715\u00d76=4290\"],\n  [\"622\u00d73=1866\", \"138\u00d79=1242\"],\n  [\"656\u00d75=3280\", \"667\u00d72=1334\"],\n  [\"413\u00d78=3304\", \"713\u00d79=6417\"],\n  [\"324\u00d78=2592\", \"662\u00d73=1986\"],\n  [\"543\u00d72=1086\", \"676\u00d75=3380\"],\n  [\"733\u00d78=5864\", \"206\u00d74=824\"],\n  [\"830\u00d75=4150\", \"807\u00d75=4035\"],\n  [\"216\u00d77=1512\", \"529\u00d74=2116\"],\n  [\"724\u00d79=6516\", \"259\u00d73=777\"],\n  [\"422\u00d75=2110\", \"104\u00d74=416\"],\n  [\"447\u00d78=3576\", \"194\u00d79=1746\"],\n  [\"175\u00d73=525\", \"359\u00d73=1077\"],\n  [\"272\u00d78=2176\", \"368\u00d75=1840\"],\n  [\"363\u00d74=1452\", \"961\u00d75=4805\"],\n  [\"616\u00d73=1848\", \"115\u00d76=690\"],\n  [\"692\u00d73=2076\", \"587\u00d74=2348\"],\n  [\"335\u00d73=1005\", \"176\u00d76=1056\"],\n  [\"386\u00d77=2702\", \"783\u00d79=7047\"],\n  [\"643\u00d73=1929\", \"448\u00d73=1344\"],\n  [\"741\u00d79=6669\", \"625\u00d72=1250\"],\n  [\"990\u00d73=2970\", \"839\u00d77=5873\"],\n  [\"663\u00d76=3978\", \"534\u00d75=2670\"],\n  [\"507\u00d74=2028\", \"981\u00d77=6867\"],\n  [\"714\u00d77=4998\", \"469\u00d79=4221\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit x one-digit multiplication problems/answers\n# in the table cells with the new values from the commit diff.\n# Each old value is unique within the document, so Find/Replace against\n# the whole document range is sufficient and keeps the original run\n# formatting intact (wdReplaceOne = 1 semantics per match via Execute).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"276\u00d77=1932\", \"715\u00d76=4290\"),\n    @(\"622\u00d73=1866\", \"138\u00d79=1242\"),\n    @(\"656\u00d75=3280\", \"667\u00d72=1334\"),\n    @(\"413\u00d78=3304\", \"713\u00d79=6417\"),\n    @(\"324\u00d78=2592\", \"662\u00d73=1986\"),\n    @(\"543\u00d72=1086\", \"676\u00d75=3380\"),\n    @(\"733\u00d78=5864\", \"206\u00d74=824\"),\n    @(\"830\u00d75=4150\", \"807\u00d75=4035\"),\n    @(\"216\u00d77=1512\", \"529\u00d74=2116\"),\n    @(\"724\u00d79=6516\", \"259\u00d73=777\"),\n    @(\"422\u00d75=2110\", \"104\u00d74=416\"),\n    @(\"447\u00d78=3576\", \"194\u00d79=1746\"),\n    @(\"175\u00d73=525\", \"359\u00d73=1077\"),\n    @(\"272\u00d78=2176\", \"368\u00d75=1840\"),\n    @(\"363\u00d74=1452\", \"961\u00d75=4805\"),\n    @(\"616\u00d73=1848\", \"115\u00d76=690\"),\n    @(\"692\u00d73=2076\", \"587\u00d74=2348\"),\n    @(\"335\u00d73=1005\", \"176\u00d76=1056\"),\n    @(\"386\u00d77=2702\", \"783\u00d79=7047\"),\n    @(\"643\u00d73=1929\", \"448\u00d73=1344\"),\n    @(\"741\u00d79=6669\", \"625\u00d72=1250\"),\n    @(\"990\u00d73=2970\", \"839\u00d77=5873\"),\n    @(\"663\u00d76=3978\", \"534\u00d75=2670\"),\n    @(\"507\u00d74=2028\", \"981\u00d77=6867\"),\n    @(\"714\u00d77=4998\", \"469\u00d79=4221\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $result = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $result) {\n        Write-Output \"WARNING: replacement not found for '$old'\"\n    }\n}\n"}
